# Update the directive names in the certificate template.
$d = $word.ActiveDocument

# 1) "{#certs}Dña. " -> "{#certs}D. "
$d.Content.Find.Execute("{#certs}Dña. ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{#certs}D. ", 2)

# 2) "María Teresa Juan Díaz, secretaria" -> "José Manuel Sanz Molinero, secretario"
$d.Content.Find.Execute("María Teresa Juan Díaz, secretaria", $true, $false, $false, $false, $false,
                         $true, 1, $false, "José Manuel Sanz Molinero, secretario", 2)

# 3) ", con " + "N.I.F." + " " -> ", con N.I.F. " (collapse three runs into the text of the first)
$d.Content.Find.Execute("N.I.F.  ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "N.I.F. ", 2)
